$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data for columns I (I0) and J (IF), keyed by row number
$data = @"
2,7,7
3,6,7
4,8,8
5,7,7
6,7,8
7,8,8
8,8,8
9,8,8
10,8,8
11,8,8
12,8,8
13,8,9
14,8,8
15,8,8
16,8,8
17,8,8
18,6,6
19,8,8
20,8,9
21,6,7
22,6,7
23,6,6
24,9,9
25,9,9
26,8,8
27,7,8
28,8,8
29,7,8
30,7,7
31,7,7
32,7,7
33,7,7
34,8,8
35,7,7
36,7,8
37,7,7
38,6,7
39,9,9
40,7,7
41,6,6
42,6,6
43,7,7
44,11,11
45,8,8
46,8,8
47,8,8
48,7,7
49,6,6
50,4,5
51,7,8
52,6,6
53,5,5
54,4,5
55,8,8
56,9,9
57,6,6
58,9,9
59,9,9
60,8,8
61,8,8
62,6,7
63,9,9
64,8,8
65,8,8
66,6,6
67,5,5
68,7,7
69,4,4
"@

foreach ($line in $data -split "`n") {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $iVal = [int]$parts[1]
    $jVal = [int]$parts[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}

Write-Host "Applied I0/IF columns to $($data -split "`n" | Where-Object { $_.Trim() -ne "" } | Measure-Object | Select-Object -ExpandProperty Count) rows"
